# Add new columns I (I0) and J (IF) to Sheet1, mirroring the header style
# of the existing H column ("IP") and filling in the per-row numeric data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy H1's format (bold, centered, bordered) into
#     I1 and J1, then overwrite the copied text with the new header labels.
$h1 = $ws.Range("H1")
$i1 = $ws.Range("I1")
$j1 = $ws.Range("J1")

$h1.Copy($i1)
$h1.Copy($j1)

$i1.Value = "I0"
$j1.Value = "IF"

# --- Data rows 2-42: column I then column J values.
$data = @(
  @(2,8,9),
  @(3,1,3),
  @(4,10,10),
  @(5,10,11),
  @(6,7,7),
  @(7,9,9),
  @(8,8,8),
  @(9,7,7),
  @(10,9,9),
  @(11,8,9),
  @(12,11,11),
  @(13,9,9),
  @(14,6,7),
  @(15,7,8),
  @(16,8,8),
  @(17,8,9),
  @(18,4,7),
  @(19,8,8),
  @(20,9,9),
  @(21,9,9),
  @(22,9,9),
  @(23,9,9),
  @(24,7,7),
  @(25,11,11),
  @(26,8,9),
  @(27,7,9),
  @(28,8,8),
  @(29,11,11),
  @(30,10,10),
  @(31,7,7),
  @(32,6,8),
  @(33,9,9),
  @(34,6,7),
  @(35,10,11),
  @(36,7,7),
  @(37,9,9),
  @(38,4,6),
  @(39,9,9),
  @(40,8,8),
  @(41,7,9),
  @(42,1,2)
)

foreach ($entry in $data) {
  $row = $entry[0]
  $iVal = $entry[1]
  $jVal = $entry[2]
  $ws.Cells.Item($row, 9).Value = $iVal
  $ws.Cells.Item($row, 10).Value = $jVal
}
